$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.944.51'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.553.33'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.36'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.08%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.546'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.88%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.71'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0823'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('E12').Value = '  +5.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.63'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.78%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.944.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.542.87'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('E16').Value = '  +4.69%  '
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.046.87'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.67'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.37%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0988'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.97'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.76%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '255.32'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('E24').Value = '  +1.31%  '
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '27.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.99'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.27%  '
$ws.Range('E30').Value = '  -1.17%  '
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '158.99'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.73%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.16'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.75'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0805'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('E36').Value = '  -1.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.90'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +12.55%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '25.96'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +10.33%  '
$ws.Range('E39').Value = '  -1.47%  '
$ws.Range('E40').Value = '  -0.52%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.12'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +34.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.45'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.22%  '
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.101.35'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.26%  '
$ws.Range('E45').Value = '  -2.51%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '86.93'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('E48').Value = '  +3.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.801.62'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.95'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '103.79'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.77%  '
